$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Tauros" row of data that used to live at row 5
$ws.Range("A5:J5").ClearContents()

# Row 6 (Rasmusmon) gains a formatted (but empty) L6 cell
$ws.Range("L6").HorizontalAlignment = -4131

# New, mostly-empty row 7 with just a formatted L7 cell
$ws.Range("L7").HorizontalAlignment = -4131

# Re-add the Tauros entry much further down, at row 129, with an extra
# duplicated type column and a move-list entry
$ws.Range("A129").Value = 128
$ws.Range("B129").Value = "Tauros"
$ws.Range("C129").Value = "NORMAL"
$ws.Range("D129").Value = "NORMAL"
$ws.Range("E129").Value = 75
$ws.Range("F129").Value = 100
$ws.Range("G129").Value = 95
$ws.Range("H129").Value = 110
$ws.Range("I129").Value = 101
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = "Tackle,SwordsDance"
$ws.Range("L129").HorizontalAlignment = -4131

# Update the view: zoom level and active selection
$win = $excel.ActiveWindow
$win.Zoom = 70
[void]$ws.Range("K113").Select()
